$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 350
$ws.Range("J7").Value = 350
$ws.Range("L7").Value = 350
$ws.Range("N7").Value = -574
$ws.Range("H14").Value = 350
$ws.Range("J14").Value = 350
$ws.Range("L14").Value = 350
$ws.Range("N14").Value = -732
$ws.Range("H32").Value = 2781.7144
$ws.Range("I32").Value = 2928.3333
$ws.Range("J32").Value = 2671.75
$ws.Range("K32").Value = 2928.3333
$ws.Range("L32").Value = 2671.75
$ws.Range("M32").Value = -2602.3333
$ws.Range("N32").Value = -3323.75
$ws.Range("H74").Value = 7734.846
$ws.Range("I74").Value = 6150.6
$ws.Range("K74").Value = 6150.6
$ws.Range("M74").Value = -5214.6
$ws.Range("H77").Value = 7734.846
$ws.Range("I77").Value = 6150.6
$ws.Range("K77").Value = 30753
$ws.Range("M77").Value = -26073
$ws.Range("H86").Value = 3034.4666
$ws.Range("I86").Value = 1986
$ws.Range("J86").Value = 3951.875
$ws.Range("K86").Value = 1986
$ws.Range("L86").Value = 3951.875
$ws.Range("M86").Value = -863
$ws.Range("N86").Value = -6197.875
$ws.Range("H89").Value = 3034.4666
$ws.Range("I89").Value = 1986
$ws.Range("J89").Value = 3951.875
$ws.Range("K89").Value = 9930
$ws.Range("L89").Value = 19759.375
$ws.Range("M89").Value = -4314
$ws.Range("N89").Value = -30991.375
$ws.Range("H113").Value = 3817.575
$ws.Range("I113").Value = 5685.9375
$ws.Range("J113").Value = 2572
$ws.Range("K113").Value = 5685.9375
$ws.Range("L113").Value = 2572
$ws.Range("M113").Value = -2431.9375
$ws.Range("N113").Value = -9080
$ws.Range("H125").Value = 2306.353
$ws.Range("I125").Value = 2450
$ws.Range("J125").Value = 2275.5715
$ws.Range("K125").Value = 22050
$ws.Range("L125").Value = 20480.1435
$ws.Range("M125").Value = -19590
$ws.Range("N125").Value = -25400.1435
$ws.Range("H127").Value = 1799.6666
$ws.Range("I127").Value = 1156.7142
$ws.Range("K127").Value = 3470.1426
$ws.Range("M127").Value = 1489.8574
$ws.Range("H132").Value = 68993.375
$ws.Range("I132").Value = 72568.55
$ws.Range("J132").Value = 10003
$ws.Range("K132").Value = 217705.65
$ws.Range("L132").Value = 30009
$ws.Range("M132").Value = -215175.65
$ws.Range("N132").Value = -35069
$ws.Range("H137").Value = 2721.1853
$ws.Range("J137").Value = 3142.6843
$ws.Range("L137").Value = 9428.052899999999
$ws.Range("N137").Value = -14528.0529

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5959597.5
$ws.Range("I32").Value = 6673029
$ws.Range("J32").Value = 14332.833
$ws.Range("K32").Value = 6673029
$ws.Range("L32").Value = 14332.833
$ws.Range("M32").Value = -6672742
$ws.Range("N32").Value = -14906.833
$ws.Range("H61").Value = 8523.143
$ws.Range("I61").Value = 8404.200000000001
$ws.Range("J61").Value = 8647.041999999999
$ws.Range("K61").Value = 8404.200000000001
$ws.Range("L61").Value = 8647.041999999999
$ws.Range("M61").Value = -8192.200000000001
$ws.Range("N61").Value = -9071.041999999999
$ws.Range("H74").Value = 3516.25
$ws.Range("I74").Value = 2992.0652
$ws.Range("K74").Value = 2992.0652
$ws.Range("M74").Value = -2118.0652
$ws.Range("H77").Value = 3516.25
$ws.Range("I77").Value = 2992.0652
$ws.Range("K77").Value = 14960.326
$ws.Range("M77").Value = -10592.326
$ws.Range("H132").Value = 6442.373
$ws.Range("I132").Value = 5455.6665
$ws.Range("K132").Value = 16366.9995
$ws.Range("M132").Value = -13836.9995
$ws.Range("H136").Value = 8523.143
$ws.Range("I136").Value = 8404.200000000001
$ws.Range("J136").Value = 8647.041999999999
$ws.Range("K136").Value = 25212.6
$ws.Range("L136").Value = 25941.126
$ws.Range("M136").Value = -22662.6
$ws.Range("N136").Value = -31041.126

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 45333.332
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 45333.332
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 45333.332
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -45963.332
$ws.Range("H79").Value = 45333.332
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 45333.332
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 45333.332
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -47517.332
$ws.Range("H99").Value = 1748.3462
$ws.Range("I99").Value = 1351.9166
$ws.Range("J99").Value = 6505.5
$ws.Range("K99").Value = 1351.9166
$ws.Range("L99").Value = 6505.5
$ws.Range("M99").Value = 146.0834
$ws.Range("N99").Value = -9501.5
$ws.Range("H134").Value = 2915.8918
$ws.Range("I134").Value = 2058.1091
$ws.Range("K134").Value = 6174.327300000001
$ws.Range("M134").Value = -3639.327300000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7301.7236
$ws.Range("I31").Value = 1764.1666
$ws.Range("J31").Value = 10738.827
$ws.Range("K31").Value = 1764.1666
$ws.Range("L31").Value = 10738.827
$ws.Range("M31").Value = -1469.1666
$ws.Range("N31").Value = -11328.827
$ws.Range("H34").Value = 7301.7236
$ws.Range("I34").Value = 1764.1666
$ws.Range("J34").Value = 10738.827
$ws.Range("K34").Value = 1764.1666
$ws.Range("L34").Value = 10738.827
$ws.Range("M34").Value = -1562.1666
$ws.Range("N34").Value = -11142.827
$ws.Range("H99").Value = 7659.222
$ws.Range("I99").Value = 8572.166999999999
$ws.Range("K99").Value = 8572.166999999999
$ws.Range("M99").Value = -7074.166999999999
$ws.Range("H126").Value = 7659.222
$ws.Range("I126").Value = 8572.166999999999
$ws.Range("K126").Value = 25716.501
$ws.Range("M126").Value = -23246.501

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 3000
$ws.Range("J54").Value = 3000
$ws.Range("L54").Value = 9000
$ws.Range("N54").Value = -10118
$ws.Range("H117").Value = 4583.263
$ws.Range("I117").Value = 1447
$ws.Range("J117").Value = 7405.9
$ws.Range("K117").Value = 4341
$ws.Range("L117").Value = 22217.7
$ws.Range("M117").Value = -899
$ws.Range("N117").Value = -29101.7

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4600
$ws.Range("I70").Value = 5466.6665
$ws.Range("J70").Value = 2000
$ws.Range("K70").Value = 5466.6665
$ws.Range("L70").Value = 2000
$ws.Range("M70").Value = -5196.6665
$ws.Range("N70").Value = -2540
$ws.Range("H73").Value = 4600
$ws.Range("I73").Value = 5466.6665
$ws.Range("J73").Value = 2000
$ws.Range("K73").Value = 5466.6665
$ws.Range("L73").Value = 2000
$ws.Range("M73").Value = -4530.6665
$ws.Range("N73").Value = -3872
$ws.Range("H97").Value = 1024.2963
$ws.Range("I97").Value = 616.6923
$ws.Range("K97").Value = 616.6923
$ws.Range("M97").Value = -120.6923
$ws.Range("H107").Value = 416.46667
$ws.Range("I107").Value = 271.16666
$ws.Range("J107").Value = 997.6667
$ws.Range("K107").Value = 271.16666
$ws.Range("L107").Value = 997.6667
$ws.Range("M107").Value = 1648.83334
$ws.Range("N107").Value = -4837.6667
$ws.Range("H121").Value = 5174.75
$ws.Range("J121").Value = 5174.75
$ws.Range("L121").Value = 5174.75
$ws.Range("N121").Value = -8668.75
$ws.Range("H122").Value = 4470.1665
$ws.Range("I122").Value = 2229.4119
$ws.Range("K122").Value = 6688.2357
$ws.Range("M122").Value = -4238.2357
$ws.Range("H132").Value = 2418.6155
$ws.Range("I132").Value = 2371
$ws.Range("J132").Value = 2525.75
$ws.Range("K132").Value = 7113
$ws.Range("L132").Value = 7577.25
$ws.Range("M132").Value = -4583
$ws.Range("N132").Value = -12637.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 6865.3335
$ws.Range("I16").Value = 6865.3335
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 6865.3335
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -6695.3335
$ws.Range("N16").ClearContents()
$ws.Range("H22").Value = 58612.668
$ws.Range("I22").Value = 250687.5
$ws.Range("J22").Value = 3734.1428
$ws.Range("K22").Value = 250687.5
$ws.Range("L22").Value = 3734.1428
$ws.Range("M22").Value = -250392.5
$ws.Range("N22").Value = -4324.1428
$ws.Range("H25").Value = 70000
$ws.Range("I25").Value = 70000
$ws.Range("K25").Value = 70000
$ws.Range("M25").Value = -69770
$ws.Range("H27").Value = 58612.668
$ws.Range("I27").Value = 250687.5
$ws.Range("J27").Value = 3734.1428
$ws.Range("K27").Value = 250687.5
$ws.Range("L27").Value = 3734.1428
$ws.Range("M27").Value = -250580.5
$ws.Range("N27").Value = -3948.1428
$ws.Range("H35").Value = 1766
$ws.Range("I35").Value = 1611.75
$ws.Range("J35").Value = 3000
$ws.Range("K35").Value = 1611.75
$ws.Range("L35").Value = 3000
$ws.Range("M35").Value = -1275.75
$ws.Range("N35").Value = -3672
$ws.Range("H38").Value = 15663.333
$ws.Range("I38").Value = 15000
$ws.Range("K38").Value = 15000
$ws.Range("M38").Value = -14590
$ws.Range("H100").Value = 6129.769
$ws.Range("J100").Value = 12702.1
$ws.Range("L100").Value = 12702.1
$ws.Range("N100").Value = -13784.1
$ws.Range("H122").Value = 3500.24
$ws.Range("I122").Value = 3235.1904
$ws.Range("J122").Value = 4891.75
$ws.Range("K122").Value = 9705.5712
$ws.Range("L122").Value = 14675.25
$ws.Range("M122").Value = -7255.5712
$ws.Range("N122").Value = -19575.25
$ws.Range("H132").Value = 4993.3613
$ws.Range("J132").Value = 6253.6924
$ws.Range("L132").Value = 18761.0772
$ws.Range("N132").Value = -23821.0772

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 36666.5
$ws.Range("I8").Value = 70000
$ws.Range("K8").Value = 70000
$ws.Range("M8").Value = -69860
$ws.Range("H122").Value = 3590.9473
$ws.Range("I122").Value = 3326.9375
$ws.Range("J122").Value = 4999
$ws.Range("K122").Value = 9980.8125
$ws.Range("L122").Value = 14997
$ws.Range("M122").Value = -7530.8125
$ws.Range("N122").Value = -19897
$ws.Range("H132").Value = 13710612
$ws.Range("I132").Value = 5191.375
$ws.Range("J132").Value = 23244818
$ws.Range("K132").Value = 15574.125
$ws.Range("L132").Value = 69734454
$ws.Range("M132").Value = -13044.125
$ws.Range("N132").Value = -69739514
$ws.Range("H136").Value = 7220988
$ws.Range("I136").Value = 7702099.5
$ws.Range("J136").Value = 4316.3335
$ws.Range("K136").Value = 23106298.5
$ws.Range("L136").Value = 12949.0005
$ws.Range("M136").Value = -23103748.5
$ws.Range("N136").Value = -18049.0005
